$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column O (15) width, to host the "countifs" header/results column
$ws.Columns.Item(15).ColumnWidth = 19.95

# Header cell P19: reuse the same look as the neighbouring headers (M19:O19)
$ws.Range("O19").Copy() | Out-Null
$ws.Range("P19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("P19").Value = "countifs"

# New COUNTIFS formulas in P20:P22
$ws.Range("P20").Formula = '=COUNTIFS(N10:N16,"kharid sakte",F10:F16,"lg")'
$ws.Range("P21").Formula = '=COUNTIFS(N10:N16,"nhi kharid sakte",F10:F16,"apple")'
$ws.Range("P22").Formula = '=COUNTIFS(N11:N17,"koshish krte",F11:F17,"samsung")'

# Update the visible selection to match the newly added header cells
$ws.Range("O19:P19").Select() | Out-Null
